$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.379.03"
$ws.Range("E2").Value = "  +4.96%  "
$ws.Range("D3").Value = "1.813.64"
$ws.Range("E3").Value = "  +5.66%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "344.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.10%  "
$ws.Range("E7").Value = "  +3.28%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3489"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.44%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "48.86"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.91%  "
$ws.Range("E10").Value = "  +3.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07725"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.46%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.003"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.08"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +9.98%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.616"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.00%  "
$ws.Range("D15").Value = "1.814.96"
$ws.Range("E15").Value = "  +5.84%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.221"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +4.05%  "
$ws.Range("E17").Value = "  +3.77%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06733"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.47%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "85.95"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.88%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.001"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.14%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.61"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.49%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.584"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +8.20%  "
$ws.Range("E23").Value = "  +1.63%  "
$ws.Range("D24").Value = "27.389.89"
$ws.Range("E24").Value = "  +5.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.463"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.40%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.666"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +8.70%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "21.96"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +14.26%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.466"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +9.54%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "154.50"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("D30").Value = "2.019.13"
$ws.Range("E30").Value = "  +5.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "135.76"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.311"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.35%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.14%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "13.91"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.08769"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.693"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.22%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.615"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.92%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6973"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +13.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2272"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +6.26%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02404"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +4.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.06473"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +3.93%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.931"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +4.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.306"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.29%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "14.70"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6523"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.001"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.042"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.45%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.178"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +7.83%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "132.57"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.07326"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.85%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.43"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.37%  "
